$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.809.39'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '2.618.50'
$ws.Range('E3').Value = '  -3.29%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '589.41'
$ws.Range('E5').Value = '  -2.73%  '
$ws.Range('D6').Value = '163.97'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  -3.72%  '
$ws.Range('D9').Value = '2.618.79'
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('D10').Value = '0.142'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').Value = '0.359'
$ws.Range('D13').Value = '5.20'
$ws.Range('E13').Value = '  -1.69%  '
$ws.Range('D14').Value = '27.31'
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('D15').Value = '3.117.29'
$ws.Range('D16').Value = '0.0000180'
$ws.Range('E16').Value = '  -3.50%  '
$ws.Range('D17').Value = '66.940.20'
$ws.Range('E17').Value = '  -1.69%  '
$ws.Range('D18').Value = '2.618.93'
$ws.Range('E18').Value = '  -3.11%  '
$ws.Range('D19').Value = '11.90'
$ws.Range('E19').Value = '  +0.70%  '
$ws.Range('D20').Value = '7.99'
$ws.Range('E20').Value = '  +5.16%  '
$ws.Range('D21').Value = '356.29'
$ws.Range('E21').Value = '  -3.73%  '
$ws.Range('D22').Value = '4.32'
$ws.Range('E22').Value = '  -3.84%  '
$ws.Range('D23').Value = '4.66'
$ws.Range('E23').Value = '  -5.80%  '
$ws.Range('D24').Value = '10.87'
$ws.Range('E24').Value = '  +7.47%  '
$ws.Range('D25').Value = '1.94'
$ws.Range('E25').Value = '  -6.68%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '70.59'
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '0.0000100'
$ws.Range('E30').Value = '  -3.42%  '
$ws.Range('D31').Value = '546.86'
$ws.Range('E31').Value = '  -5.15%  '
$ws.Range('D32').Value = '7.88'
$ws.Range('E32').Value = '  -3.38%  '
$ws.Range('D33').Value = '1.35'
$ws.Range('E33').Value = '  -4.54%  '
$ws.Range('E34').Value = '  -4.96%  '
$ws.Range('E35').Value = '  +1.40%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -4.86%  '
$ws.Range('D38').Value = '157.32'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').Value = '19.09'
$ws.Range('E39').Value = '  -3.78%  '
$ws.Range('D40').Value = '0.365'
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('D41').Value = '5.18'
$ws.Range('E41').Value = '  -3.93%  '
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -4.81%  '
$ws.Range('D43').Value = '17.89'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '2.45'
$ws.Range('E45').Value = '  -5.97%  '
$ws.Range('D46').Value = '40.22'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').Value = '0.0₆0295'
$ws.Range('E47').Value = '  -4.55%  '
$ws.Range('D48').Value = '0.583'
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('D49').Value = '151.15'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('D50').Value = '3.79'
$ws.Range('E50').Value = '  -2.71%  '
$ws.Range('E51').Value = '  -3.38%  '
